$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 4 data rows (old rows 2-5) so that the data that used to
# start at row 6 becomes row 2, shifting everything else up accordingly.
$ws.Range("A2:E5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
